# Insert a new data row at row 256 (pushing the existing rows 256-309 down
# to 257-310), then populate the new row with its data. This mirrors the
# source diff, which shows every row from 256 through 309 taking on the
# values that used to belong to the row above it, a brand-new row appears
# at (what becomes) row 256, and the sheet's dimension grows from R309 to
# R310.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(256).Insert()

$ws.Range("A256").Value = 10
$ws.Range("B256").Value = "Vega Modelo de Temuco"
$ws.Range("C256").Value = "La Araucanía"
$ws.Range("D256").Value = 44694
$ws.Range("E256").Value = 9
$ws.Range("F256").Value = 100112044
$ws.Range("G256").Value = "Perejil"
$ws.Range("H256").Value = "Sin especificar"
$ws.Range("I256").Value = "Primera"
$ws.Range("J256").Value = 65
$ws.Range("K256").Value = 4000
$ws.Range("L256").Value = 4000
$ws.Range("M256").Value = 4000
$ws.Range("N256").Value = "$/docena de atados (3 kilos)"
$ws.Range("O256").Value = "Provincia de Cautín"
$ws.Range("P256").Value = 1333
$ws.Range("Q256").Value = 3
$ws.Range("R256").Value = "Hortaliza"
